# Updates the cryptos list price (D) and 1h volume change (E) columns
# Values are plain text (not numbers/percentages), so NumberFormat is forced
# to Text ('@') before assignment and reset to the default 'Normal' style
# afterwards so no residual number-format/style change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Address, $Text) {
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = '@'
    $range.Value = $Text
    $range.Style = 'Normal'
}

Set-TextValue $ws 'D2' '43.694.29'
Set-TextValue $ws 'D3' '2.245.81'
Set-TextValue $ws 'E3' '  +0.72%  '
Set-TextValue $ws 'E4' '  +0.06%  '
Set-TextValue $ws 'D5' '322.71'
Set-TextValue $ws 'E5' '  +1.94%  '
Set-TextValue $ws 'D6' '101.24'
Set-TextValue $ws 'E6' '  +0.42%  '
Set-TextValue $ws 'D7' '0.579'
Set-TextValue $ws 'E7' '  -1.16%  '
Set-TextValue $ws 'E8' '  +0.05%  '
Set-TextValue $ws 'D9' '0.555'
Set-TextValue $ws 'E9' '  -0.92%  '
Set-TextValue $ws 'E10' '  +1.37%  '
Set-TextValue $ws 'E11' '  +1.08%  '
Set-TextValue $ws 'D12' '7.70'
Set-TextValue $ws 'E12' '  +0.74%  '
Set-TextValue $ws 'E13' '  -2.12%  '
Set-TextValue $ws 'D14' '2.587.83'
Set-TextValue $ws 'E14' '  +0.76%  '
Set-TextValue $ws 'D15' '0.858'
Set-TextValue $ws 'E15' '  -0.14%  '
Set-TextValue $ws 'D16' '14.23'
Set-TextValue $ws 'E16' '  -0.53%  '
Set-TextValue $ws 'D17' '2.242.50'
Set-TextValue $ws 'E17' '  +0.30%  '
Set-TextValue $ws 'D18' '43.614.90'
Set-TextValue $ws 'E18' '  +1.43%  '
Set-TextValue $ws 'D19' '13.66'
Set-TextValue $ws 'E19' '  -5.16%  '
Set-TextValue $ws 'E20' '  +2.65%  '
Set-TextValue $ws 'D21' '6.56'
Set-TextValue $ws 'E21' '  +1.01%  '
Set-TextValue $ws 'D22' '65.23'
Set-TextValue $ws 'E22' '  -0.17%  '
Set-TextValue $ws 'D23' '3.16'
Set-TextValue $ws 'E23' '  -0.74%  '
Set-TextValue $ws 'D24' '236.57'
Set-TextValue $ws 'E24' '  -0.60%  '
Set-TextValue $ws 'E25' '  +1.94%  '
Set-TextValue $ws 'E26' '  -0.04%  '
Set-TextValue $ws 'D27' '10.13'
Set-TextValue $ws 'E27' '  +1.40%  '
Set-TextValue $ws 'E28' '  -1.93%  '
Set-TextValue $ws 'D29' '36.92'
Set-TextValue $ws 'E29' '  +7.04%  '
Set-TextValue $ws 'E30' '  -1.07%  '
Set-TextValue $ws 'D31' '160.41'
Set-TextValue $ws 'E31' '  +4.14%  '
Set-TextValue $ws 'D32' '20.17'
Set-TextValue $ws 'E32' '  -1.44%  '
Set-TextValue $ws 'E33' '  -1.68%  '
Set-TextValue $ws 'E34' '  -2.96%  '
Set-TextValue $ws 'D35' '3.15'
Set-TextValue $ws 'E35' '  -1.11%  '
Set-TextValue $ws 'E36' '  +8.76%  '
Set-TextValue $ws 'E37' '  -0.17%  '
Set-TextValue $ws 'D38' '0.119'
Set-TextValue $ws 'E38' '  -1.64%  '
Set-TextValue $ws 'E39' '  +2.19%  '
Set-TextValue $ws 'D40' '4.25'
Set-TextValue $ws 'E40' '  -4.02%  '
Set-TextValue $ws 'D41' '15.61'
Set-TextValue $ws 'E41' '  +18.66%  '
Set-TextValue $ws 'E42' '  -1.44%  '
Set-TextValue $ws 'E43' '  +0.29%  '
Set-TextValue $ws 'D44' '1.804.62'
Set-TextValue $ws 'E44' '  +0.14%  '
Set-TextValue $ws 'E45' '  -2.69%  '
Set-TextValue $ws 'D46' '82.43'
Set-TextValue $ws 'E46' '  -6.24%  '
Set-TextValue $ws 'E47' '  +5.91%  '
Set-TextValue $ws 'D48' '74.52'
Set-TextValue $ws 'E48' '  -3.08%  '
Set-TextValue $ws 'D49' '5.20'
Set-TextValue $ws 'E49' '  -2.48%  '
Set-TextValue $ws 'D50' '58.65'
Set-TextValue $ws 'E50' '  -0.71%  '
Set-TextValue $ws 'D51' '103.21'
Set-TextValue $ws 'E51' '  -0.02%  '
